# Updated symbol list on Sat Dec 31 11:28:03 UTC 2022 with GitHub Actions
#
# The workbook stores numeric-looking prices as literal text (inline
# strings) rather than real numbers. To preserve that text typing when
# writing through the Excel object model (which would otherwise infer a
# Number from a numeric-looking string and also round-trip trailing
# zeros away), briefly mark the cell as Text before assigning the value,
# then restore General/Normal formatting so no visible style change is
# left behind.
function Set-TextValue {
    param($ws, $cellRef, $val)
    $r = $ws.Range($cellRef)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.NumberFormat = "General"
    $r.Style = "Normal"
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Price (column D) updates ---
Set-TextValue $ws "D2"  "246.06"
Set-TextValue $ws "D3"  "26.15"
Set-TextValue $ws "D4"  "5.099"
Set-TextValue $ws "D5"  "0.05610"
Set-TextValue $ws "D6"  "6.481"
Set-TextValue $ws "D7"  "3.022"
Set-TextValue $ws "D8"  "0.8117"
Set-TextValue $ws "D9"  "0.8474"
Set-TextValue $ws "D11" "0.03214"
Set-TextValue $ws "D13" "0.09405"
Set-TextValue $ws "D14" "0.001510"

# Row 15 (One / ONE): "Worst in 24h" badge removed from the Volume label
$ws.Range("E15").Value = "14OneONE"

Set-TextValue $ws "D16" "0.006110"
Set-TextValue $ws "D17" "3.557"
Set-TextValue $ws "D19" "0.3183"
Set-TextValue $ws "D20" "0.06974"
Set-TextValue $ws "D22" "3.743"
Set-TextValue $ws "D23" "0.04691"
Set-TextValue $ws "D25" "0.001247"
Set-TextValue $ws "D27" "0.00009599"
Set-TextValue $ws "D41" "0.1353"

# Rows 42/43 (CEJI / KickToken) swapped places in the ranking
$ws.Range("B42").Value = "KickToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
Set-TextValue $ws "D42" "0.006109"
$ws.Range("E42").Value = "41KickTokenKICK"

$ws.Range("B43").Value = "CEJI"
$ws.Range("C43").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
Set-TextValue $ws "D43" "0.002660"
$ws.Range("E43").Value = "42CEJICEJI"

Set-TextValue $ws "D44" "0.008630"
Set-TextValue $ws "D45" "0.00005290"
Set-TextValue $ws "D47" "0.1330"

# Row 47 (CoinbaseStockToken / COIN): "Worst in 24h" badge added to the Volume label
$ws.Range("E47").Value = "46CoinbaseStockTokenCOINWorstin24h"
